$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the Validator_Agent "Instructions" cell (D7) -------------------
# The cell contains rich text: a plain run, a bold run around "{final_response}",
# and a trailing plain run. We append a new sentence to the end of that
# trailing run while preserving the existing rich-text formatting.

$cell = $ws.Range("D7")
$originalText = $cell.Characters().Text()

$appendText = "`nIf you find the response to not be asking the user for input, such as if there are no questions for the user, automatically send the 'Advisor_Agent' a response of 'Ok' on behalf of the user. "

$newText = $originalText + $appendText
$cell.Characters().Text = $newText

# Re-apply the bold formatting to the " {final_response} " run that is lost
# when the whole cell text is reassigned above.
$boldToken = "{final_response}"
$tokenIdx0 = $originalText.IndexOf($boldToken)
$boldStart1 = $tokenIdx0 # 1-based start is tokenIdx0 - 1 + 1 = tokenIdx0
$boldLen = $boldToken.Length + 2

$boldRange = $cell.Characters($boldStart1, $boldLen)
$boldRange.Font.Bold = $true
$boldRange.Font.Size = 11
$boldRange.Font.Name = "Calibri"

# Re-apply the (non-bold) formatting of the trailing run, which now also
# contains the newly appended sentence.
$restStart1 = $boldStart1 + $boldLen
$totalLen = $newText.Length
$restLen = $totalLen - $restStart1 + 1

$restRange = $cell.Characters($restStart1, $restLen)
$restRange.Font.Bold = $false
$restRange.Font.Size = 11
$restRange.Font.Name = "Calibri"

# --- Update the view state to match the edited cell -------------------------
$sheetView = $ws.Application.ActiveWindow
$ws.Range("D7").Select()
$excel.ActiveWindow.ScrollRow = 6
